$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1 "NAME" is unchanged) ---
$ws.Range("B1").Value = "EMAIL"
$ws.Range("C1").Value = "LEVEL"

# --- Update data rows with the new student records ---
$ws.Range("A2").Value = "Marcio"
$ws.Range("B2").Value = "marcioc424@gmail.com"
$ws.Range("C2").Value = "master"

$ws.Range("A3").Value = "Teste"
$ws.Range("B3").Value = "teste@gmail.com"
$ws.Range("C3").Value = "starter"

$ws.Range("A4").Value = "Teste 2"
$ws.Range("B4").Value = "teste2@gmail.com"
$ws.Range("C4").Value = "intermediate"

# --- Add hyperlinks on the email cells (mailto links) ---
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:marcioc424@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:teste@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:teste2@gmail.com") | Out-Null

# --- Re-apply the original uniform cell formatting (single Consolas style)
#     to every cell, since the extra 10pt Consolas style is no longer used
#     and the hyperlink cells should keep the same look as the rest ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1:C4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Drop the now-unused "Hyperlink" named cell style that Excel
#     auto-creates when a hyperlink is added ---
try {
    $wb.Styles.Item("Hyperlink").Delete() | Out-Null
} catch {
}

$excel.CutCopyMode = 0

# --- Update the active selection shown in the sheet view ---
$ws.Range("C9").Select() | Out-Null

Write-Host "done"
